$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1068.1067
$ws.Range("I15").Value = 1068.1067
$ws.Range("K15").Value = 3204.3201
$ws.Range("M15").Value = -3035.3201

$ws.Range("H103").Value = 620.1429000000001
$ws.Range("I103").Value = 321.33334
$ws.Range("J103").Value = 739.6667
$ws.Range("K103").Value = 964.0000200000001
$ws.Range("L103").Value = 2219.0001
$ws.Range("M103").Value = -378.0000200000001
$ws.Range("N103").Value = -3391.0001

$ws.Range("H113").Value = 67910220
$ws.Range("I113").Value = 27779628
$ws.Range("J113").Value = 100014700
$ws.Range("K113").Value = 27779628
$ws.Range("L113").Value = 100014700
$ws.Range("M113").Value = -27776374
$ws.Range("N113").Value = -100021208

$ws.Range("H121").Value = 5385.4287
$ws.Range("J121").Value = 5385.4287
$ws.Range("L121").Value = 16156.2861
$ws.Range("N121").Value = -19650.2861

$ws.Range("N131").ClearContents()
$ws.Range("H131").Value = 2137.5
$ws.Range("I131").Value = 2137.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6412.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1372.5

$ws.Range("H132").Value = 1867.0197
$ws.Range("I132").Value = 1837.9375
$ws.Range("K132").Value = 5513.8125
$ws.Range("M132").Value = -2983.8125

$ws.Range("H137").Value = 3143.3618
$ws.Range("I137").Value = 3402.625
$ws.Range("J137").Value = 2872.8262
$ws.Range("K137").Value = 10207.875
$ws.Range("L137").Value = 8618.4786
$ws.Range("M137").Value = -7657.875
$ws.Range("N137").Value = -13718.4786

$ws.Range("H138").Value = 2710033.5
$ws.Range("I138").Value = 4258
$ws.Range("J138").Value = 3854784.5
$ws.Range("K138").Value = 12774
$ws.Range("L138").Value = 11564353.5
$ws.Range("M138").Value = -7634
$ws.Range("N138").Value = -11574633.5

$ws.Range("H141").Value = 2390.4644
$ws.Range("I141").Value = 1723.2609
$ws.Range("K141").Value = 5169.7827
$ws.Range("M141").Value = 10.21730000000025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1789454.5
$ws.Range("I32").Value = 1926061.5
$ws.Range("K32").Value = 1926061.5
$ws.Range("M32").Value = -1925774.5

$ws.Range("H45").Value = 5757.077
$ws.Range("I45").Value = 2426.8572
$ws.Range("J45").Value = 9642.333000000001
$ws.Range("K45").Value = 2426.8572
$ws.Range("L45").Value = 9642.333000000001
$ws.Range("M45").Value = -2049.8572
$ws.Range("N45").Value = -10396.333

$ws.Range("H122").Value = 11585.593
$ws.Range("I122").Value = 12502.652
$ws.Range("J122").Value = 6312.5
$ws.Range("K122").Value = 37507.956
$ws.Range("L122").Value = 18937.5
$ws.Range("M122").Value = -35057.956
$ws.Range("N122").Value = -23837.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4062.5417
$ws.Range("I94").Value = 2176.923
$ws.Range("J94").Value = 6291
$ws.Range("K94").Value = 2176.923
$ws.Range("L94").Value = 6291
$ws.Range("M94").Value = -1725.923
$ws.Range("N94").Value = -7193

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5913.36
$ws.Range("I31").Value = 2621.8845
$ws.Range("J31").Value = 9479.125
$ws.Range("K31").Value = 2621.8845
$ws.Range("L31").Value = 9479.125
$ws.Range("M31").Value = -2326.8845
$ws.Range("N31").Value = -10069.125

$ws.Range("H34").Value = 5913.36
$ws.Range("I34").Value = 2621.8845
$ws.Range("J34").Value = 9479.125
$ws.Range("K34").Value = 2621.8845
$ws.Range("L34").Value = 9479.125
$ws.Range("M34").Value = -2419.8845
$ws.Range("N34").Value = -9883.125

$ws.Range("H122").Value = 1718.3214
$ws.Range("I122").Value = 1413.0476
$ws.Range("K122").Value = 4239.142800000001
$ws.Range("M122").Value = -1789.142800000001

$ws.Range("H132").Value = 7031.6816
$ws.Range("I132").Value = 4091.3333
$ws.Range("J132").Value = 9067.308000000001
$ws.Range("K132").Value = 12273.9999
$ws.Range("L132").Value = 27201.924
$ws.Range("M132").Value = -9743.999899999999
$ws.Range("N132").Value = -32261.924

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4902.643
$ws.Range("I68").Value = 2880.6667
$ws.Range("J68").Value = 5454.091
$ws.Range("K68").Value = 8642.000100000001
$ws.Range("L68").Value = 16362.273
$ws.Range("M68").Value = -7831.000100000001
$ws.Range("N68").Value = -17984.273

$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 2500
$ws.Range("L69").Value = 7500
$ws.Range("N69").Value = -9122

$ws.Range("M70").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0

$ws.Range("H71").Value = 4902.643
$ws.Range("I71").Value = 2880.6667
$ws.Range("J71").Value = 5454.091
$ws.Range("K71").Value = 25926.0003
$ws.Range("L71").Value = 49086.819
$ws.Range("M71").Value = -21870.0003
$ws.Range("N71").Value = -57198.819

$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 2500
$ws.Range("L72").Value = 22500
$ws.Range("N72").Value = -30612

$ws.Range("M73").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0

$ws.Range("H113").Value = 2388.9167
$ws.Range("I113").Value = 1330.8334
$ws.Range("J113").Value = 2917.9583
$ws.Range("K113").Value = 3992.5002
$ws.Range("L113").Value = 8753.874899999999
$ws.Range("M113").Value = -1822.5002
$ws.Range("N113").Value = -13093.8749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5390.3076
$ws.Range("I132").Value = 3154.0588
$ws.Range("K132").Value = 9462.1764
$ws.Range("M132").Value = -6932.1764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4148.5
$ws.Range("I93").Value = 3982.8572
$ws.Range("K93").Value = 3982.8572
$ws.Range("M93").Value = -2734.8572

$ws.Range("H122").Value = 2284.0408
$ws.Range("I122").Value = 1690.5385
$ws.Range("J122").Value = 4598.7
$ws.Range("K122").Value = 5071.6155
$ws.Range("L122").Value = 13796.1
$ws.Range("M122").Value = -2621.6155
$ws.Range("N122").Value = -18696.1

$ws.Range("H132").Value = 17247866
$ws.Range("I132").Value = 38464510
$ws.Range("J132").Value = 9340.25
$ws.Range("K132").Value = 115393530
$ws.Range("L132").Value = 28020.75
$ws.Range("M132").Value = -115391000
$ws.Range("N132").Value = -33080.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 77006000
$ws.Range("I136").Value = 166669580
$ws.Range("K136").Value = 500008740
$ws.Range("M136").Value = -500006190
